# Daily attendance processing - 2025-11-11 19:42:51
# Swap the order of names in the "Recorded By" (column G) cells so that
# "dnasr281@gmail.com" is listed first when it appears alongside exactly
# one other recorder (e.g. "System" or "admin@admin.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -eq 2 -and ($parts[0] -eq $target -or $parts[1] -eq $target) -and $parts[0] -ne $parts[1]) {
            $newValue = $parts[1] + ", " + $parts[0]
            $cell.Value2 = $newValue
        }
    }
}
